$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 22

$ws.Cells.Item($row, 1).Value = "'2024-01-05"
$ws.Cells.Item($row, 2).Value = "18:33:32"
$ws.Cells.Item($row, 3).Value = "Friday"
$ws.Cells.Item($row, 4).Value = "'00"
$ws.Cells.Item($row, 5).Value = 140629
$ws.Cells.Item($row, 6).Value = 142886
$ws.Cells.Item($row, 7).Value = 172376
$ws.Cells.Item($row, 8).Value = 147197
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 118148
$ws.Cells.Item($row, 11).Value = 224466
$ws.Cells.Item($row, 12).Value = 248844
$ws.Cells.Item($row, 13).Value = 184715
$ws.Cells.Item($row, 14).Value = 110132
$ws.Cells.Item($row, 15).Value = 40490
$ws.Cells.Item($row, 16).Value = 30813
$ws.Cells.Item($row, 17).Value = 72404
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 41707
$ws.Cells.Item($row, 20).Value = -1
